$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D28").Value = "2021년 결산"
$ws.Range("E28").Value = "https://ropiens.tistory.com/174"

$ws.Range("D45").Value = "상관계수"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/96"

$ws.Range("D46").Value = "[씨젠] 2021년 01월, 생물정보학(Bioinformatics 채용), Bioinformatics Engineer"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/425"
